# Loan RBI, Variable Instalments
# On the "Repayment schedule" sheet, a new column is inserted before column N
# ("In Advance"), shifting the existing "In Advance", "Date" and "Outstanding"
# columns one place to the right. The new column inherits the width of the
# column to its left ("Original").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

$originalColWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $originalColWidth

# The "Repayment schedule" sheet becomes the active/selected sheet and tab,
# with cell R7 selected.
$ws.Activate()
$ws.Range("R7").Select() | Out-Null
